$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 16
$ws1.Range("F4").Value = 968

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 16
$ws4.Range("F4").Value = 968
